$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

$ws.Range("B6").Formula = "=208.1-171.1"
$ws.Range("B9").Formula = "=143.3-183"
$ws.Range("B16").Formula = "=63.7-50.5"
$ws.Range("B21").Formula = "=41.3-48.2"

$ws.Range("B22").Select()
